$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the two new sheets, positioned after RunModes ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "LoginTest"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "LoginTest2"

# --- LoginTest (sheet2) content ---
$ws2.Range("A1").Value = "Fname"
$ws2.Range("B1").Value = "Lname"
$ws2.Range("C1").Value = "Email"
$ws2.Range("D1").Value = "CrePswd"
$ws2.Range("E1").Value = "ConPswd"
$ws2.Range("A1:E1").Interior.Color = $ws1.Range("A1").Interior.Color

$ws2.Range("A2").Value = "Narayana1"
$ws2.Range("B2").Value = "M1"
$ws2.Range("C2").Value = "n1@n.com"
$ws2.Range("D2").Value = "*****1"
$ws2.Range("E2").Value = "*****1"

$ws2.Range("A3").Value = "Narayana2"
$ws2.Range("B3").Value = "M2"
$ws2.Range("C3").Value = "n1@n.com"
$ws2.Range("D3").Value = "*****2"
$ws2.Range("E3").Value = "*****2"

$ws2.Range("A4").Value = "Narayana3"
$ws2.Range("B4").Value = "M3"
$ws2.Range("C4").Value = "n1@n.com"
$ws2.Range("D4").Value = "*****3"
$ws2.Range("E4").Value = "*****3"

$ws2.Range("A5").Value = "Narayana4"
$ws2.Range("B5").Value = "M4"
$ws2.Range("C5").Value = "n1@n.com"
$ws2.Range("D5").Value = "*****4"
$ws2.Range("E5").Value = "*****4"

$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:n1@n.com")
$ws2.Hyperlinks.Add($ws2.Range("C3:C5"), "mailto:n1@n.com", "", "", "n1@n.com")
$ws2.Range("C2:C5").Style = "Hyperlink"

$ws2.Columns.Item(1).ColumnWidth = 10.71

[void]$ws2.Range("E1").Select()

# --- LoginTest2 (sheet3) stays empty ---
[void]$ws3.Range("F4").Select()

# --- sheet1 selection update ---
[void]$ws1.Range("C12").Select()

# --- Make LoginTest2 the active tab ---
[void]$ws3.Activate()
